# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which contain duplicate/mirrored data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    4  = 12518
    5  = 1278
    7  = 32
    10 = 199
    11 = 453
    12 = 60
    16 = 377
    17 = 4172
    19 = 13
    20 = 942
    21 = 21
    22 = 125
    23 = 64
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
